$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("C1").Value = 0.06903616931152196
$ws.Range("D1").Value = 0
$ws.Range("E1").Value = 214.1613132505006

# Row 2
$ws.Range("B2").Value = 0.9
$ws.Range("C2").Value = 0.0360770973831241
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 187.356844242832
